# paises.xlsx refresh: "Update countries & provincias Spain"
# Re-ranked country rows (shared-string / row reshuffle) + updated case counts,
# as of "Datos actualizados a 31 de Marzo de 2020 a las 20:20".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country column (A) reshuffles: rows whose country changed rank/position ---
$ws.Cells.Item(60, 1).Value = 'Irak'
$ws.Cells.Item(61, 1).Value = 'Catar'
$ws.Cells.Item(62, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(63, 1).Value = 'Egipto'
$ws.Cells.Item(64, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(98, 1).Value = 'Uzbekistan'
$ws.Cells.Item(99, 1).Value = 'Malta'
$ws.Cells.Item(100, 1).Value = 'Islas Feroe'
$ws.Cells.Item(101, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(128, 1).Value = 'Monaco'
$ws.Cells.Item(129, 1).Value = 'Banglades'
$ws.Cells.Item(131, 1).Value = 'Uganda'
$ws.Cells.Item(132, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(133, 1).Value = 'Macao'
$ws.Cells.Item(134, 1).Value = 'Puerto Rico'
$ws.Cells.Item(135, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(136, 1).Value = 'Jamaica'
$ws.Cells.Item(137, 1).Value = 'Guatemala'
$ws.Cells.Item(138, 1).Value = 'Zambia'
$ws.Cells.Item(139, 1).Value = 'Barbados'
$ws.Cells.Item(140, 1).Value = 'Togo'
$ws.Cells.Item(181, 1).Value = 'Republica del Chad'
$ws.Cells.Item(182, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(187, 1).Value = 'Cabo Verde'
$ws.Cells.Item(188, 1).Value = 'San Bartolome'
$ws.Cells.Item(190, 1).Value = 'Montserrat'
$ws.Cells.Item(191, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(192, 1).Value = 'Fiyi'
$ws.Cells.Item(193, 1).Value = 'Somalia'
$ws.Cells.Item(194, 1).Value = 'Nepal'
$ws.Cells.Item(195, 1).Value = 'Butan'
$ws.Cells.Item(196, 1).Value = 'Gambia'
$ws.Cells.Item(197, 1).Value = 'Nicaragua'
$ws.Cells.Item(202, 1).Value = 'Islas Virgenes Britanicas'

# --- Updated statistics: B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes ---
$ws.Range('B9').Value = 52128
$ws.Range('C9').Value = 7578
$ws.Range('D9').Value = 9444
$ws.Range('E9').Value = 39161
$ws.Range('F9').Value = 5565
$ws.Range('G9').Value = 499
$ws.Range('H9').Value = 3523

$ws.Range('B16').Value = 10122
$ws.Range('C16').Value = 504
$ws.Range('D16').Value = 1095
$ws.Range('E16').Value = 8899
$ws.Range('F16').Value = 198
$ws.Range('G16').Value = 20
$ws.Range('H16').Value = 128

$ws.Range('B20').Value = 4831
$ws.Range('C20').Value = 136
$ws.Range('D20').Value = 163
$ws.Range('E20').Value = 4648
$ws.Range('F20').Value = 83
$ws.Range('G20').Value = 4
$ws.Range('H20').Value = 20

$ws.Range('B25').Value = 3257
$ws.Range('C25').Value = 256
$ws.Range('D25').Value = 45
$ws.Range('E25').Value = 3181
$ws.Range('F25').Value = 64
$ws.Range('G25').Value = 8
$ws.Range('H25').Value = 31

$ws.Range('B60').Value = 694
$ws.Range('C60').Value = 64
$ws.Range('D60').Value = 170
$ws.Range('E60').Value = 474
$ws.Range('F60').Value = 0
$ws.Range('G60').Value = 4
$ws.Range('H60').Value = 50

$ws.Range('B61').Value = 693
$ws.Range('C61').Value = 0
$ws.Range('D61').Value = 51
$ws.Range('E61').Value = 641
$ws.Range('F61').Value = 6
$ws.Range('G61').Value = 0
$ws.Range('H61').Value = 1

$ws.Range('B62').Value = 664
$ws.Range('C62').Value = 53
$ws.Range('D62').Value = 61
$ws.Range('E62').Value = 597
$ws.Range('F62').Value = 2
$ws.Range('G62').Value = 1
$ws.Range('H62').Value = 6

$ws.Range('B63').Value = 656
$ws.Range('C63').Value = 0
$ws.Range('D63').Value = 150
$ws.Range('E63').Value = 465
$ws.Range('F63').Value = 0
$ws.Range('G63').Value = 0
$ws.Range('H63').Value = 41

$ws.Range('B64').Value = 647
$ws.Range('C64').Value = 58
$ws.Range('D64').Value = 74
$ws.Range('E64').Value = 572
$ws.Range('F64').Value = 2
$ws.Range('G64').Value = 0
$ws.Range('H64').Value = 1

$ws.Range('B98').Value = 172
$ws.Range('C98').Value = 23
$ws.Range('D98').Value = 7
$ws.Range('E98').Value = 163
$ws.Range('F98').Value = 8
$ws.Range('G98').Value = 0
$ws.Range('H98').Value = 2

$ws.Range('B99').Value = 169
$ws.Range('C99').Value = 13
$ws.Range('D99').Value = 2
$ws.Range('E99').Value = 167
$ws.Range('F99').Value = 2
$ws.Range('G99').Value = 0
$ws.Range('H99').Value = 0

$ws.Range('B100').Value = 169
$ws.Range('C100').Value = 1
$ws.Range('D100').Value = 74
$ws.Range('E100').Value = 95
$ws.Range('F100').Value = 3
$ws.Range('G100').Value = 0
$ws.Range('H100').Value = 0

$ws.Range('B101').Value = 168
$ws.Range('C101').Value = 0
$ws.Range('D101').Value = 6
$ws.Range('E101').Value = 161
$ws.Range('F101').Value = 0
$ws.Range('G101').Value = 0
$ws.Range('H101').Value = 1

$ws.Range('B105').Value = 143
$ws.Range('C105').Value = 21
$ws.Range('D105').Value = 17
$ws.Range('E105').Value = 124
$ws.Range('F105').Value = 5
$ws.Range('G105').Value = 0
$ws.Range('H105').Value = 2

$ws.Range('B128').Value = 52
$ws.Range('C128').Value = 3
$ws.Range('D128').Value = 2
$ws.Range('E128').Value = 49
$ws.Range('F128').Value = 0
$ws.Range('G128').Value = 0
$ws.Range('H128').Value = 1

$ws.Range('B129').Value = 51
$ws.Range('C129').Value = 2
$ws.Range('D129').Value = 25
$ws.Range('E129').Value = 21
$ws.Range('F129').Value = 1
$ws.Range('G129').Value = 0
$ws.Range('H129').Value = 5

$ws.Range('B131').Value = 44
$ws.Range('C131').Value = 11
$ws.Range('D131').Value = 0
$ws.Range('E131').Value = 44
$ws.Range('F131').Value = 0
$ws.Range('G131').Value = 0
$ws.Range('H131').Value = 0

$ws.Range('B132').Value = 43
$ws.Range('C132').Value = 0
$ws.Range('D132').Value = 6
$ws.Range('E132').Value = 37
$ws.Range('F132').Value = 0
$ws.Range('G132').Value = 0
$ws.Range('H132').Value = 0

$ws.Range('B133').Value = 41
$ws.Range('C133').Value = 3
$ws.Range('D133').Value = 10
$ws.Range('E133').Value = 31
$ws.Range('F133').Value = 0
$ws.Range('G133').Value = 0
$ws.Range('H133').Value = 0

$ws.Range('B134').Value = 39
$ws.Range('C134').Value = 0
$ws.Range('D134').Value = 1
$ws.Range('E134').Value = 36
$ws.Range('F134').Value = 0
$ws.Range('G134').Value = 0
$ws.Range('H134').Value = 2

$ws.Range('B135').Value = 36
$ws.Range('C135').Value = 0
$ws.Range('D135').Value = 0
$ws.Range('E135').Value = 36
$ws.Range('F135').Value = 2
$ws.Range('G135').Value = 0
$ws.Range('H135').Value = 0

$ws.Range('B136').Value = 36
$ws.Range('C136').Value = 0
$ws.Range('D136').Value = 2
$ws.Range('E136').Value = 33
$ws.Range('F136').Value = 0
$ws.Range('G136').Value = 0
$ws.Range('H136').Value = 1

$ws.Range('B137').Value = 36
$ws.Range('C137').Value = 0
$ws.Range('D137').Value = 10
$ws.Range('E137').Value = 25
$ws.Range('F137').Value = 1
$ws.Range('G137').Value = 0
$ws.Range('H137').Value = 1

$ws.Range('B138').Value = 35
$ws.Range('C138').Value = 0
$ws.Range('D138').Value = 0
$ws.Range('E138').Value = 35
$ws.Range('F138').Value = 0
$ws.Range('G138').Value = 0
$ws.Range('H138').Value = 0

$ws.Range('B139').Value = 34
$ws.Range('C139').Value = 1
$ws.Range('D139').Value = 0
$ws.Range('E139').Value = 34
$ws.Range('F139').Value = 0
$ws.Range('G139').Value = 0
$ws.Range('H139').Value = 0

$ws.Range('B140').Value = 34
$ws.Range('C140').Value = 4
$ws.Range('D140').Value = 10
$ws.Range('E140').Value = 23
$ws.Range('F140').Value = 0
$ws.Range('G140').Value = 0
$ws.Range('H140').Value = 1

$ws.Range('B181').Value = 7
$ws.Range('C181').Value = 2
$ws.Range('D181').Value = 0
$ws.Range('E181').Value = 7
$ws.Range('F181').Value = 0
$ws.Range('G181').Value = 0
$ws.Range('H181').Value = 0

$ws.Range('B182').Value = 7
$ws.Range('C182').Value = 0
$ws.Range('D182').Value = 0
$ws.Range('E182').Value = 7
$ws.Range('F182').Value = 0
$ws.Range('G182').Value = 0
$ws.Range('H182').Value = 0

$ws.Range('B187').Value = 6
$ws.Range('C187').Value = 0
$ws.Range('D187').Value = 0
$ws.Range('E187').Value = 5
$ws.Range('F187').Value = 0
$ws.Range('G187').Value = 0
$ws.Range('H187').Value = 1

$ws.Range('B188').Value = 6
$ws.Range('C188').Value = 0
$ws.Range('D188').Value = 1
$ws.Range('E188').Value = 5
$ws.Range('F188').Value = 0
$ws.Range('G188').Value = 0
$ws.Range('H188').Value = 0

$ws.Range('B193').Value = 5
$ws.Range('C193').Value = 2
$ws.Range('D193').Value = 1
$ws.Range('E193').Value = 4
$ws.Range('F193').Value = 0
$ws.Range('G193').Value = 0
$ws.Range('H193').Value = 0

$ws.Range('B194').Value = 5
$ws.Range('C194').Value = 0
$ws.Range('D194').Value = 1
$ws.Range('E194').Value = 4
$ws.Range('F194').Value = 0
$ws.Range('G194').Value = 0
$ws.Range('H194').Value = 0

$ws.Range('B195').Value = 4
$ws.Range('C195').Value = 0
$ws.Range('D195').Value = 0
$ws.Range('E195').Value = 4
$ws.Range('F195').Value = 0
$ws.Range('G195').Value = 0
$ws.Range('H195').Value = 0

$ws.Range('B197').Value = 4
$ws.Range('C197').Value = 0
$ws.Range('D197').Value = 0
$ws.Range('E197').Value = 3
$ws.Range('F197').Value = 0
$ws.Range('G197').Value = 0
$ws.Range('H197').Value = 1

$ws.Range('B202').Value = 3
$ws.Range('C202').Value = 1
$ws.Range('D202').Value = 0
$ws.Range('E202').Value = 3
$ws.Range('F202').Value = 0
$ws.Range('G202').Value = 0
$ws.Range('H202').Value = 0

# --- Footer timestamp ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 31 de Marzo de 2020 a las 20:20'
